$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q1" right before the "总计" sheet
# ---------------------------------------------------------------------
$refSheet   = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# match the page margins used by the other quarter sheets (0.75/1/0.5 in)
$newSheet.PageSetup.LeftMargin   = 54
$newSheet.PageSetup.RightMargin  = 54
$newSheet.PageSetup.TopMargin    = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# copy the header style (bold + box border) from an existing quarter sheet
$refSheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# copy the index-column style (bold + box border) down for the 8 data rows
$refSheet.Range("A2:A9").Copy($newSheet.Range("A2:A9"))

$fundRows = @(
    @("003886", "汇安丰利灵活配置混合A",               "7.39", "94.04", "3.60", "0.2660", 9),
    @("003887", "汇安丰利灵活配置混合C",               "4.33", "94.04", "3.60", "0.1559", 9),
    @("011410", "中信建投量化进取6个月持有期混合A",   "9.13", "93.80", "1.00", "0.0913", 6),
    @("011411", "中信建投量化进取6个月持有期混合C",   "2.15", "93.80", "1.00", "0.0215", 6),
    @("004194", "招商中证1000指数增强A",               "1.76", "94.40", "1.04", "0.0183", 10),
    @("004195", "招商中证1000指数增强C",               "0.68", "94.40", "1.04", "0.0071", 10),
    @("003854", "汇安丰华灵活配置混合A",               "0.19", "29.16", "1.56", "0.0030", 8),
    @("003855", "汇安丰华灵活配置混合C",               "0.19", "29.16", "1.56", "0.0030", 8)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2

    # columns B..G are stored as TEXT in this workbook (even the numeric
    # looking ones), so force text via NumberFormat "@" then restore the
    # default "Normal" style so no stray number-format sticks to the cell.
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Style = "Normal"

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" summary row at the top of the "总计" sheet,
#    shifting the existing 2021-Q4 / 2021-Q3 / 2021-Q2 rows down by one
# ---------------------------------------------------------------------
$ts = $wb.Worksheets.Item("总计")

# extend the index-column style down to the new last row (row 5)
$ts.Range("A4").Copy($ts.Range("A5"))

# shift the three existing data rows down by one row
$ts.Range("B2:D4").Copy($ts.Range("B3:D5"))

# renumber the index column for the shifted rows
$ts.Range("A3").Value = 1
$ts.Range("A4").Value = 2
$ts.Range("A5").Value = 3

# write the new first data row
$ts.Range("A2").Value = 0
$ts.Range("B2").Value = "2022-Q1"
$ts.Range("C2").Value = 8
$ts.Range("D2").Value = 0.57
